# ISA.xlsx update:
#  - Rename Sheet1 -> "8-bit" and Sheet2 -> "16-bit" (planning the migration
#    to a preliminary 16-bit architecture while the 8-bit design is piped).
#  - Move the active selection on the 8-bit sheet from N10 to O5.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Name = "8-bit"
$ws2.Name = "16-bit"

$ws1.Activate() | Out-Null
$ws1.Range("O5").Select() | Out-Null
